$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit inserts one new data record for "Papa" (row 876) and pushes every
# existing record from old row 876 through old row 977 down by one row
# (to new rows 877 through 978). Columns A, B, C, E, F, G, Q, R stay constant
# for this whole block, so we only need to shift D, H, I, J, K, L, M, N, O, P.

# First, populate the brand-new last row (978) with the constant columns,
# since that row does not exist yet.
$ws.Cells.Item(978, 1).Value2  = $ws.Cells.Item(977, 1).Value2   # A - Mercado ID
$ws.Cells.Item(978, 2).Value2  = $ws.Cells.Item(977, 2).Value2   # B - Mercado
$ws.Cells.Item(978, 3).Value2  = $ws.Cells.Item(977, 3).Value2   # C - Region
$ws.Cells.Item(978, 5).Value2  = $ws.Cells.Item(977, 5).Value2   # E - Codreg
$ws.Cells.Item(978, 6).Value2  = $ws.Cells.Item(977, 6).Value2   # F - Categoria ID
$ws.Cells.Item(978, 7).Value2  = $ws.Cells.Item(977, 7).Value2   # G - Categoria
$ws.Cells.Item(978, 17).Value2 = $ws.Cells.Item(977, 17).Value2  # Q - Kg o Unidades
$ws.Cells.Item(978, 18).Value2 = $ws.Cells.Item(977, 18).Value2  # R - Clasificacion

# Shift rows 876..977 down into 877..978, working bottom-up so we never
# clobber a source row before it has been copied.
for ($i = 978; $i -ge 877; $i--) {
    $src = $i - 1
    $ws.Cells.Item($i, 4).Value2  = $ws.Cells.Item($src, 4).Value2   # D - Fecha
    $ws.Cells.Item($i, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat
    $ws.Cells.Item($i, 8).Value2  = $ws.Cells.Item($src, 8).Value2   # H - Variedad
    $ws.Cells.Item($i, 9).Value2  = $ws.Cells.Item($src, 9).Value2   # I - Calidad
    $ws.Cells.Item($i, 10).Value2 = $ws.Cells.Item($src, 10).Value2  # J - Volumen
    $ws.Cells.Item($i, 11).Value2 = $ws.Cells.Item($src, 11).Value2  # K - Precio minimo
    $ws.Cells.Item($i, 12).Value2 = $ws.Cells.Item($src, 12).Value2  # L - Precio maximo
    $ws.Cells.Item($i, 13).Value2 = $ws.Cells.Item($src, 13).Value2  # M - Precio promedio ponderado
    $ws.Cells.Item($i, 14).Value2 = $ws.Cells.Item($src, 14).Value2  # N - Unidad de comercializacion
    $ws.Cells.Item($i, 15).Value2 = $ws.Cells.Item($src, 15).Value2  # O - Origen
    $ws.Cells.Item($i, 16).Value2 = $ws.Cells.Item($src, 16).Value2  # P - Precio $/Kg
}

# Finally, overwrite row 876 with the new record's data.
$ws.Cells.Item(876, 4).Value2  = 44918
$ws.Cells.Item(876, 8).Value2  = "Rosara"
$ws.Cells.Item(876, 9).Value2  = "1a nueva(o)"
$ws.Cells.Item(876, 10).Value2 = 350
$ws.Cells.Item(876, 11).Value2 = 13000
$ws.Cells.Item(876, 12).Value2 = 13000
$ws.Cells.Item(876, 13).Value2 = 13000
$ws.Cells.Item(876, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(876, 16).Value2 = 520
